# "dadata added to address" - apply the textual edits described by the diff.
#
# Replace-All (wdReplaceAll) is document-wide in this engine regardless of
# which Range/Paragraph object invoked Find, so several short target
# strings (e.g. "12", "111111", ".56") would otherwise also hit unrelated
# numbers elsewhere in the document (zip codes, registration numbers...).
# To avoid that collateral damage every replacement below is done with a
# single Find/Replace (wdReplaceOne) against an explicit $d.Range(start,end)
# scoped tightly to the paragraph holding the text. The start/end offsets
# are re-read from the live Paragraphs collection immediately before each
# call (rather than hard-coded), since earlier replacements shift the
# character offsets of everything that follows them.

$d = $word.ActiveDocument

function Replace-InParagraph($index, $old, $new) {
    $p = $d.Paragraphs.Item($index)
    $start = $p.Range.Start
    $end = $p.Range.End
    $r = $d.Range($start, $end)
    $r.Find.Execute($old, $true, $false, $false, $false, $false, `
                     $true, 0, $false, $new, 1) | Out-Null
}

# Paragraph 8: software/product title (first occurrence)
Replace-InParagraph 8 `
    "Программное обеспечение для автоматического формирования документа и электронного документооборота" `
    "Система формирования документов и электронного документооборота"

# Paragraph 9: application number "12___...___" -> "23___...___"
Replace-InParagraph 9 "12______________________________________________" "23______________________________________________"

# Paragraph 13: address - postal code "111111" -> "413111"
Replace-InParagraph 13 "111111 " "413111 "

# Paragraph 13: address - apartment no. & phone "56, +79898887716" -> "23, 89271180894"
Replace-InParagraph 13 "56, +79898887716" "23, 89271180894"

# Paragraph 14: passport / identity document details
Replace-InParagraph 14 `
    "1111 № 123453 выдан МВД России по Саратовской области в городе Энгельсе Дата выдачи: 22.07.2022 " `
    "1111 № 111122 выдан отделом уфмс росии по саратовской области в городе Энгельсе Дата выдачи: 12.01.2022 "

# Paragraph 22: date "07.06.2024" -> "10.06.2024"
Replace-InParagraph 22 "07.06.2024" "10.06.2024"

# Paragraph 35: application number in the second table "12" -> "23"
Replace-InParagraph 35 "12" "23"

# Paragraph 38: software/product title (second occurrence, closing guillemet)
Replace-InParagraph 38 `
    "Программное обеспечение для автоматического формирования документа и электронного документооборота»" `
    "Система формирования документов и электронного документооборота»"

# Paragraph 50: address block - postal code "111111" -> "413111"
Replace-InParagraph 50 "Россия, 111111, " "Россия, 413111, "

# Paragraph 50: address block - apartment no. ".56" -> ".23"
Replace-InParagraph 50 ".56" ".23"
